# Settings back to qfin22 presentation - recalibration
$wb = $excel.ActiveWorkbook

$wsAR      = $wb.Worksheets.Item("AR")
$wsSETAR   = $wb.Worksheets.Item("SETAR")
$wsGARCH   = $wb.Worksheets.Item("GARCH")
$wsTARCH   = $wb.Worksheets.Item("TARCH")
$wsARTARCH = $wb.Worksheets.Item("AR_TARCH")

# ---------------- AR ----------------
$wsAR.Range("B2").Value = "-1.997186433010218E-05"
$wsAR.Range("B2").Style = "Normal"
$wsAR.Range("B3").Value = 0.780208627016614
$wsAR.Range("B4").Value = 0.0001159205106355226
$wsAR.Range("B5").Value = "[0.9999999999999999, 0.20228377255522503, 0.16267650023453922, 0.12726831164056143, 0.12518901139632232, 0.3671427030140748, 0.14554010513854237, 0.13596012033354346, 0.11368610256417404, 0.09561022053693262, 0.06325720027941474, 0.07060911345884024, 0.07237921058027443, 0.10966245609834162, 0.1401772100435194, 0.10041024811725854, 0.0668958304710907, 0.02784687768461108, 0.05558976733911963, 0.08990897194346244]"

# ---------------- SETAR ----------------
$wsSETAR.Range("B4").Value = -0.004096972279127004
$wsSETAR.Range("B5").Value = 0.6917522748811119
$wsSETAR.Range("B6").Value = "8.583817529437564E-05"
$wsSETAR.Range("B6").Style = "Normal"
$wsSETAR.Range("B7").Value = 0.004595438951464681
$wsSETAR.Range("B8").Value = 0.6250544330132463
$wsSETAR.Range("B9").Value = "6.343841559805757E-05"
$wsSETAR.Range("B9").Style = "Normal"
$wsSETAR.Range("B10").Value = "[0.9999999999999999, 0.22705239661064963, 0.17988360403380674, 0.18370218496126733, 0.15960375513752784, 0.252271328419774, 0.15539288602018475, 0.13370812630222562, 0.10512209370261381, 0.1191410545107181, 0.09481054235663657, 0.08325822052800276, 0.07161164162383542, 0.09637738743992473, 0.11788421245582827, 0.0826186371668529, 0.07285431965565738, 0.0576836779388001, 0.05734716468632466, 0.060289080468721674]"

# ---------------- GARCH ----------------
$wsGARCH.Range("B2").Value = -0.0001103970106513943
$wsGARCH.Range("B3").Value = "2.604993169738602E-06"
$wsGARCH.Range("B3").Style = "Normal"
$wsGARCH.Range("B4").Value = 0.09999999853467051
$wsGARCH.Range("B5").Value = 0.8799999871051004
$wsGARCH.Range("B6").Value = "[1.0, -0.0014563473723728747, -0.04086175134264834, -0.06870855205424087, -0.08130075073603629, 0.2716789365814485, -0.0002351373855570668, 0.038728912709676454, -0.0007993117053095164, -0.03483796237355698, -0.06388040172046135, -0.06413151182190131, -0.027580243464521988, 0.012534511062392874, 0.05275055167028347, 0.020187774941007177, -0.03562475749954847, -0.07968226944793948, -0.04784775323159723, 0.009156156248482598]"

# ---------------- TARCH ----------------
$wsTARCH.Range("B2").Value = -0.0001307439151577172
$wsTARCH.Range("B3").Value = "2.604993169973577E-06"
$wsTARCH.Range("B3").Style = "Normal"
$wsTARCH.Range("B4").Value = 0.09999999854369067
$wsTARCH.Range("B5").Value = 0.8749999872572933
$wsTARCH.Range("B6").Value = "[1.0, -0.004846527266885321, -0.0438405635776766, -0.07089783419345613, -0.08384042387291622, 0.2689830223481254, -0.0021539451244373393, 0.03729194196357781, -0.002052334463596128, -0.03575196879338397, -0.06425724073770196, -0.06466549597048658, -0.028214982947172115, 0.012263021079593914, 0.052496762119064455, 0.019900909227098376, -0.03596125543706273, -0.0797005179390171, -0.048048464686244755, 0.00915563909697059]"
$wsTARCH.Range("B7").Value = 0.009999985291275684

# ---------------- AR_TARCH ----------------
$wsARTARCH.Range("B2").Value = "-5.412080586049639E-05"
$wsARTARCH.Range("B2").Style = "Normal"
$wsARTARCH.Range("B3").Value = "2.371234768483483E-06"
$wsARTARCH.Range("B3").Style = "Normal"
$wsARTARCH.Range("B4").Value = 0.1000409206915611
$wsARTARCH.Range("B5").Value = 0.875004045550012
$wsARTARCH.Range("B6").Value = "[1.0, 0.021678171310174053, -0.03716825257054828, -0.06714239477651848, -0.04776216483649707, 0.22160801330035756, -0.025528235766241913, 0.019932531567171924, -0.015775755516603776, -0.030261227728658125, -0.06247205905942555, -0.058088403684193185, -0.023370608569804528, -0.0007132141744936319, 0.041379753862258964, 0.0017087979928617424, -0.03662244172896323, -0.06416890112673505, -0.04661948818274307, 0.012769301446215435]"
$wsARTARCH.Range("B7").Value = 0.01008897968694154
$wsARTARCH.Range("B9").Value = 1.235046560754107
